# Automatic update of files.
# - Bump the "Förändrad" date (column C) for every data row from 45184 to 45186.
# - Add the record id as the friendly-name second argument to every
#   HYPERLINK() formula in columns S, T, V, W, X, Y (rows that have them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 17

# Folder (under the klasma.github.io base URL) + file extension used for
# each link column.
$linkCols = @(
    @{ Col = "S"; Folder = "artfynd";         Ext = "xlsx" },
    @{ Col = "T"; Folder = "kartor";          Ext = "png"  },
    @{ Col = "V"; Folder = "klagomål";        Ext = "docx" },
    @{ Col = "W"; Folder = "klagomålsmail";   Ext = "docx" },
    @{ Col = "X"; Folder = "tillsyn";         Ext = "docx" },
    @{ Col = "Y"; Folder = "tillsynsmail";    Ext = "docx" }
)

for ($row = 2; $row -le $lastRow; $row++) {

    # Column C: bump the "changed" date from 45184 to 45186.
    $ws.Range("C$row").Value = 45186

    # Record id (e.g. "A 39255-2019") lives in column A of the same row.
    $id = $ws.Range("A$row").Value2

    if ($id) {
        foreach ($link in $linkCols) {
            $cell = $ws.Range($link.Col + $row)
            $existing = $cell.Formula

            if ($existing -and $existing -like "*HYPERLINK(*") {
                $url = "https://klasma.github.io/Logging_GOTEBORG/" + $link.Folder + "/" + $id + "." + $link.Ext
                $cell.Formula = '=HYPERLINK("' + $url + '", "' + $id + '")'
            }
        }
    }
}
